$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 6).Value2 = 1.01
$ws.Cells.Item(2, 7).Value2 = 1.02
$ws.Cells.Item(2, 8).Value2 = 230
$ws.Cells.Item(2, 9).Value2 = 1000
$ws.Cells.Item(2, 10).Value2 = 46
$ws.Cells.Item(2, 11).Value2 = 85
$ws.Cells.Item(2, 12).Value2 = 0
$ws.Cells.Item(2, 13).Value2 = 0
$ws.Cells.Item(2, 14).Value2 = 0
$ws.Cells.Item(2, 15).Value2 = 0
$ws.Cells.Item(2, 16).Value2 = 1.64
$ws.Cells.Item(2, 17).Value2 = 2.52
$ws.Cells.Item(2, 18).Value2 = 1.07
$ws.Cells.Item(2, 19).Value2 = 13.5
$ws.Cells.Item(2, 20).Value2 = 3.5
$ws.Cells.Item(2, 21).Value2 = 1.28
$ws.Cells.Item(2, 22).Value2 = 1.01
$ws.Cells.Item(2, 23).Value2 = 17
$ws.Cells.Item(2, 24).Value2 = 1000
$ws.Cells.Item(2, 25).Value2 = 1000
$ws.Cells.Item(2, 26).Value2 = 1000
$ws.Cells.Item(2, 27).Value2 = 1000
$ws.Cells.Item(2, 28).Value2 = 1000
$ws.Cells.Item(2, 29).Value2 = 1000
$ws.Cells.Item(2, 30).Value2 = 1000
$ws.Cells.Item(2, 31).Value2 = 1000
$ws.Cells.Item(2, 32).Value2 = 1.62
$ws.Cells.Item(2, 33).Value2 = 6.8
$ws.Cells.Item(2, 34).Value2 = 1000
$ws.Cells.Item(2, 35).Value2 = 1000
$ws.Cells.Item(2, 36).Value2 = 9.4
$ws.Cells.Item(2, 37).Value2 = 1000
$ws.Cells.Item(2, 38).Value2 = 650
$ws.Cells.Item(2, 39).Value2 = 1000
$ws.Cells.Item(2, 40).Value2 = 140
$ws.Cells.Item(2, 41).Value2 = 1000
$ws.Cells.Item(3, 6).Value2 = 1.61
$ws.Cells.Item(3, 7).Value2 = 1.66
$ws.Cells.Item(3, 8).Value2 = 6.8
$ws.Cells.Item(3, 9).Value2 = 7.4
$ws.Cells.Item(3, 10).Value2 = 3.95
$ws.Cells.Item(3, 11).Value2 = 4.1
$ws.Cells.Item(3, 13).Value2 = 1.08
$ws.Cells.Item(3, 14).Value2 = 3.05
$ws.Cells.Item(3, 15).Value2 = 1.46
$ws.Cells.Item(3, 16).Value2 = 1.68
$ws.Cells.Item(3, 17).Value2 = 2.32
$ws.Cells.Item(3, 18).Value2 = 1.24
$ws.Cells.Item(3, 19).Value2 = 4.9
$ws.Cells.Item(3, 20).Value2 = 2.22
$ws.Cells.Item(3, 21).Value2 = 1.68
$ws.Cells.Item(3, 22).Value2 = 1.15
$ws.Cells.Item(3, 23).Value2 = 2.42
$ws.Cells.Item(3, 25).Value2 = 20
$ws.Cells.Item(3, 26).Value2 = 85
$ws.Cells.Item(3, 27).Value2 = 390
$ws.Cells.Item(3, 28).Value2 = 6.4
$ws.Cells.Item(3, 30).Value2 = 34
$ws.Cells.Item(3, 31).Value2 = 200
$ws.Cells.Item(3, 33).Value2 = 11
$ws.Cells.Item(3, 34).Value2 = 34
$ws.Cells.Item(3, 35).Value2 = 190
$ws.Cells.Item(3, 36).Value2 = 15
$ws.Cells.Item(3, 37).Value2 = 21
$ws.Cells.Item(3, 38).Value2 = 70
$ws.Cells.Item(3, 39).Value2 = 300
$ws.Cells.Item(3, 40).Value2 = 13.5
$ws.Cells.Item(3, 41).Value2 = 370
$ws.Cells.Item(4, 6).Value2 = 2.82
$ws.Cells.Item(4, 7).Value2 = 2.86
$ws.Cells.Item(4, 8).Value2 = 2.96
$ws.Cells.Item(4, 9).Value2 = 3
$ws.Cells.Item(4, 10).Value2 = 3.2
$ws.Cells.Item(4, 12).Value2 = 1.5
$ws.Cells.Item(4, 13).Value2 = 1.1
$ws.Cells.Item(4, 14).Value2 = 3.15
$ws.Cells.Item(4, 15).Value2 = 1.43
$ws.Cells.Item(4, 16).Value2 = 1.72
$ws.Cells.Item(4, 17).Value2 = 2.3
$ws.Cells.Item(4, 18).Value2 = 1.26
$ws.Cells.Item(4, 19).Value2 = 4.4
$ws.Cells.Item(4, 20).Value2 = 1.9
$ws.Cells.Item(4, 21).Value2 = 2
$ws.Cells.Item(4, 22).Value2 = 1.5
$ws.Cells.Item(4, 23).Value2 = 1.54
$ws.Cells.Item(4, 25).Value2 = 10
$ws.Cells.Item(4, 26).Value2 = 20
$ws.Cells.Item(4, 27).Value2 = 130
$ws.Cells.Item(4, 29).Value2 = 7.4
$ws.Cells.Item(4, 30).Value2 = 13.5
$ws.Cells.Item(4, 31).Value2 = 40
$ws.Cells.Item(4, 32).Value2 = 16.5
$ws.Cells.Item(4, 33).Value2 = 12.5
$ws.Cells.Item(4, 34).Value2 = 18.5
$ws.Cells.Item(4, 35).Value2 = 85
$ws.Cells.Item(4, 36).Value2 = 44
$ws.Cells.Item(4, 37).Value2 = 34
$ws.Cells.Item(4, 38).Value2 = 55
$ws.Cells.Item(4, 39).Value2 = 130
$ws.Cells.Item(4, 40).Value2 = 85
$ws.Cells.Item(4, 41).Value2 = 600
$ws.Cells.Item(5, 6).Value2 = 1.65
$ws.Cells.Item(5, 9).Value2 = 7
$ws.Cells.Item(5, 10).Value2 = 4
$ws.Cells.Item(5, 11).Value2 = 4.2
$ws.Cells.Item(5, 14).Value2 = 3.3
$ws.Cells.Item(5, 15).Value2 = 1.41
$ws.Cells.Item(5, 16).Value2 = 1.81
$ws.Cells.Item(5, 17).Value2 = 2.18
$ws.Cells.Item(5, 18).Value2 = 1.28
$ws.Cells.Item(5, 19).Value2 = 4.2
$ws.Cells.Item(5, 20).Value2 = 2.16
$ws.Cells.Item(5, 21).Value2 = 1.84
$ws.Cells.Item(5, 22).Value2 = 1.17
$ws.Cells.Item(5, 24).Value2 = 13
$ws.Cells.Item(5, 25).Value2 = 19
$ws.Cells.Item(5, 27).Value2 = 240
$ws.Cells.Item(5, 28).Value2 = 6.8
$ws.Cells.Item(5, 29).Value2 = 9.199999999999999
$ws.Cells.Item(5, 30).Value2 = 26
$ws.Cells.Item(5, 31).Value2 = 440
$ws.Cells.Item(5, 32).Value2 = 8.6
$ws.Cells.Item(5, 33).Value2 = 9.800000000000001
$ws.Cells.Item(5, 35).Value2 = 120
$ws.Cells.Item(5, 36).Value2 = 15.5
$ws.Cells.Item(5, 37).Value2 = 19.5
$ws.Cells.Item(5, 38).Value2 = 46
$ws.Cells.Item(5, 40).Value2 = 12.5
$ws.Cells.Item(6, 6).Value2 = 1.73
$ws.Cells.Item(6, 7).Value2 = 1.76
$ws.Cells.Item(6, 9).Value2 = 6.2
$ws.Cells.Item(6, 10).Value2 = 3.85
$ws.Cells.Item(6, 11).Value2 = 4
$ws.Cells.Item(6, 12).Value2 = 1.46
$ws.Cells.Item(6, 13).Value2 = 1.08
$ws.Cells.Item(6, 14).Value2 = 3.5
$ws.Cells.Item(6, 15).Value2 = 1.36
$ws.Cells.Item(6, 16).Value2 = 1.82
$ws.Cells.Item(6, 17).Value2 = 2.08
$ws.Cells.Item(6, 18).Value2 = 1.31
$ws.Cells.Item(6, 19).Value2 = 3.9
$ws.Cells.Item(6, 20).Value2 = 1.99
$ws.Cells.Item(6, 21).Value2 = 1.82
$ws.Cells.Item(6, 22).Value2 = 1.19
$ws.Cells.Item(6, 24).Value2 = 13
$ws.Cells.Item(6, 25).Value2 = 18.5
$ws.Cells.Item(6, 26).Value2 = 55
$ws.Cells.Item(6, 27).Value2 = 180
$ws.Cells.Item(6, 30).Value2 = 24
$ws.Cells.Item(6, 31).Value2 = 110
$ws.Cells.Item(6, 32).Value2 = 9.6
$ws.Cells.Item(6, 34).Value2 = 23
$ws.Cells.Item(6, 35).Value2 = 110
$ws.Cells.Item(6, 36).Value2 = 17
$ws.Cells.Item(6, 37).Value2 = 19.5
$ws.Cells.Item(6, 38).Value2 = 44
$ws.Cells.Item(6, 39).Value2 = 160
$ws.Cells.Item(6, 40).Value2 = 12.5
$ws.Cells.Item(6, 41).Value2 = 160
$ws.Cells.Item(7, 8).Value2 = 48
$ws.Cells.Item(7, 10).Value2 = 15
$ws.Cells.Item(7, 11).Value2 = 19
$ws.Cells.Item(7, 14).Value2 = 8
$ws.Cells.Item(7, 16).Value2 = 3.5
$ws.Cells.Item(7, 17).Value2 = 1.37
$ws.Cells.Item(7, 19).Value2 = 1.9
$ws.Cells.Item(7, 20).Value2 = 3.3
$ws.Cells.Item(7, 23).Value2 = 12
$ws.Cells.Item(7, 32).Value2 = 8.199999999999999
$ws.Cells.Item(7, 33).Value2 = 980
$ws.Cells.Item(7, 36).Value2 = 8.199999999999999
$ws.Cells.Item(7, 40).Value2 = 2.4
$ws.Cells.Item(8, 6).Value2 = 1.65
$ws.Cells.Item(8, 7).Value2 = 1.67
$ws.Cells.Item(8, 8).Value2 = 7.2
$ws.Cells.Item(8, 9).Value2 = 7.4
$ws.Cells.Item(8, 12).Value2 = 1.54
$ws.Cells.Item(8, 13).Value2 = 1.11
$ws.Cells.Item(8, 14).Value2 = 2.86
$ws.Cells.Item(8, 15).Value2 = 1.51
$ws.Cells.Item(8, 16).Value2 = 1.61
$ws.Cells.Item(8, 17).Value2 = 2.56
$ws.Cells.Item(8, 18).Value2 = 1.21
$ws.Cells.Item(8, 19).Value2 = 5.1
$ws.Cells.Item(8, 20).Value2 = 2.42
$ws.Cells.Item(8, 21).Value2 = 1.61
$ws.Cells.Item(8, 22).Value2 = 1.15
$ws.Cells.Item(8, 23).Value2 = 2.5
$ws.Cells.Item(8, 27).Value2 = 340
$ws.Cells.Item(8, 28).Value2 = 5.9
$ws.Cells.Item(8, 29).Value2 = 9.4
$ws.Cells.Item(8, 30).Value2 = 30
$ws.Cells.Item(8, 31).Value2 = 160
$ws.Cells.Item(8, 32).Value2 = 8
$ws.Cells.Item(8, 34).Value2 = 34
$ws.Cells.Item(8, 35).Value2 = 180
$ws.Cells.Item(8, 36).Value2 = 16
$ws.Cells.Item(8, 37).Value2 = 23
$ws.Cells.Item(8, 38).Value2 = 65
$ws.Cells.Item(8, 39).Value2 = 280
$ws.Cells.Item(8, 40).Value2 = 17.5
$ws.Cells.Item(8, 41).Value2 = 290
$ws.Cells.Item(9, 6).Value2 = 1.52
$ws.Cells.Item(9, 7).Value2 = 1.54
$ws.Cells.Item(9, 8).Value2 = 7.6
$ws.Cells.Item(9, 10).Value2 = 4.3
$ws.Cells.Item(9, 13).Value2 = 1.07
$ws.Cells.Item(9, 15).Value2 = 1.37
$ws.Cells.Item(9, 16).Value2 = 1.91
$ws.Cells.Item(9, 17).Value2 = 2.06
$ws.Cells.Item(9, 18).Value2 = 1.33
$ws.Cells.Item(9, 20).Value2 = 2.14
$ws.Cells.Item(9, 21).Value2 = 1.76
$ws.Cells.Item(9, 22).Value2 = 1.13
$ws.Cells.Item(9, 23).Value2 = 2.8
$ws.Cells.Item(9, 24).Value2 = 14.5
$ws.Cells.Item(9, 25).Value2 = 980
$ws.Cells.Item(9, 28).Value2 = 7.2
$ws.Cells.Item(9, 30).Value2 = 1000
$ws.Cells.Item(9, 31).Value2 = 180
$ws.Cells.Item(9, 32).Value2 = 8.4
$ws.Cells.Item(9, 35).Value2 = 160
$ws.Cells.Item(9, 36).Value2 = 24
$ws.Cells.Item(9, 37).Value2 = 29
$ws.Cells.Item(9, 38).Value2 = 1000
$ws.Cells.Item(9, 40).Value2 = 9.6
$ws.Cells.Item(10, 6).Value2 = 1.46
$ws.Cells.Item(10, 7).Value2 = 1.48
$ws.Cells.Item(10, 8).Value2 = 10
$ws.Cells.Item(10, 9).Value2 = 12
$ws.Cells.Item(10, 10).Value2 = 4.4
$ws.Cells.Item(10, 11).Value2 = 4.7
$ws.Cells.Item(10, 12).Value2 = 1.42
$ws.Cells.Item(10, 13).Value2 = 1.07
$ws.Cells.Item(10, 14).Value2 = 3.6
$ws.Cells.Item(10, 15).Value2 = 1.35
$ws.Cells.Item(10, 16).Value2 = 1.92
$ws.Cells.Item(10, 17).Value2 = 2.04
$ws.Cells.Item(10, 18).Value2 = 1.36
$ws.Cells.Item(10, 19).Value2 = 3.7
$ws.Cells.Item(10, 20).Value2 = 2.28
$ws.Cells.Item(10, 21).Value2 = 1.65
$ws.Cells.Item(10, 22).Value2 = 1.1
$ws.Cells.Item(10, 23).Value2 = 3.05
$ws.Cells.Item(10, 24).Value2 = 13.5
$ws.Cells.Item(10, 25).Value2 = 28
$ws.Cells.Item(10, 26).Value2 = 90
$ws.Cells.Item(10, 27).Value2 = 430
$ws.Cells.Item(10, 28).Value2 = 7
$ws.Cells.Item(10, 29).Value2 = 10.5
$ws.Cells.Item(10, 30).Value2 = 38
$ws.Cells.Item(10, 31).Value2 = 220
$ws.Cells.Item(10, 32).Value2 = 7.6
$ws.Cells.Item(10, 33).Value2 = 10.5
$ws.Cells.Item(10, 35).Value2 = 190
$ws.Cells.Item(10, 37).Value2 = 18
$ws.Cells.Item(10, 38).Value2 = 48
$ws.Cells.Item(10, 39).Value2 = 280
$ws.Cells.Item(10, 40).Value2 = 9.199999999999999
$ws.Cells.Item(10, 41).Value2 = 300
$ws.Cells.Item(11, 6).Value2 = 3.8
$ws.Cells.Item(11, 7).Value2 = 5.2
$ws.Cells.Item(11, 15).Value2 = 1.22
$ws.Cells.Item(11, 17).Value2 = 1.63
$ws.Cells.Item(11, 25).Value2 = 1000
$ws.Cells.Item(11, 30).Value2 = 1000
$ws.Cells.Item(12, 6).Value2 = 2.74
$ws.Cells.Item(12, 7).Value2 = 2.84
$ws.Cells.Item(12, 8).Value2 = 2.96
$ws.Cells.Item(12, 10).Value2 = 3.2
$ws.Cells.Item(12, 12).Value2 = 1.6
$ws.Cells.Item(12, 13).Value2 = 1.12
$ws.Cells.Item(12, 14).Value2 = 2.68
$ws.Cells.Item(12, 15).Value2 = 1.55
$ws.Cells.Item(12, 16).Value2 = 1.55
$ws.Cells.Item(12, 17).Value2 = 2.6
$ws.Cells.Item(12, 18).Value2 = 1.19
$ws.Cells.Item(12, 19).Value2 = 5.5
$ws.Cells.Item(12, 20).Value2 = 2.06
$ws.Cells.Item(12, 21).Value2 = 1.76
$ws.Cells.Item(12, 22).Value2 = 1.49
$ws.Cells.Item(12, 23).Value2 = 1.54
$ws.Cells.Item(12, 24).Value2 = 9
$ws.Cells.Item(12, 25).Value2 = 9
$ws.Cells.Item(12, 26).Value2 = 18.5
$ws.Cells.Item(12, 28).Value2 = 8.4
$ws.Cells.Item(12, 30).Value2 = 15
$ws.Cells.Item(12, 31).Value2 = 48
$ws.Cells.Item(12, 32).Value2 = 17
$ws.Cells.Item(12, 34).Value2 = 23
$ws.Cells.Item(12, 35).Value2 = 1000
$ws.Cells.Item(12, 36).Value2 = 46
$ws.Cells.Item(12, 38).Value2 = 1000
$ws.Cells.Item(12, 39).Value2 = 180
$ws.Cells.Item(12, 40).Value2 = 50
$ws.Cells.Item(12, 41).Value2 = 65
